$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.885.43'
$ws.Range("D3").Value = '3.358.91'
$ws.Range("E3").Value = '  -2.91%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '''567.31'
$ws.Range("E5").Value = '  -2.08%  '
$ws.Range("D6").Value = '''149.25'
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("E8").Value = '  +0.27%  '
$ws.Range("D9").Value = '''7.94'
$ws.Range("E9").Value = '  +1.19%  '
$ws.Range("E10").Value = '  -1.76%  '
$ws.Range("D11").Value = '''0.413'
$ws.Range("E11").Value = '  +0.99%  '
$ws.Range("D12").Value = '3.937.10'
$ws.Range("E12").Value = '  -2.82%  '
$ws.Range("E13").Value = '  +1.09%  '
$ws.Range("D14").Value = '''28.12'
$ws.Range("E14").Value = '  -1.58%  '
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '3.357.94'
$ws.Range("E15").Value = '  -2.51%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = '''0.0000169'
$ws.Range("E16").Value = '  -1.68%  '
$ws.Range("D17").Value = '60.949.69'
$ws.Range("E17").Value = '  -3.52%  '
$ws.Range("D18").Value = '''6.38'
$ws.Range("E18").Value = '  -1.26%  '
$ws.Range("D19").Value = '''14.14'
$ws.Range("E19").Value = '  -2.70%  '
$ws.Range("E20").Value = '  -3.51%  '
$ws.Range("D21").Value = '''374.36'
$ws.Range("E21").Value = '  -4.03%  '
$ws.Range("D22").Value = '''75.44'
$ws.Range("E22").Value = '  +0.95%  '
$ws.Range("D23").Value = '''0.561'
$ws.Range("E23").Value = '  -0.41%  '
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("D25").Value = '3.506.13'
$ws.Range("E25").Value = '  -2.40%  '
$ws.Range("E26").Value = '  -5.74%  '
$ws.Range("D27").Value = '''0.177'
$ws.Range("E27").Value = '  -3.35%  '
$ws.Range("D28").Value = '''0.998'
$ws.Range("E28").Value = '  -0.27%  '
$ws.Range("E29").Value = '  -3.83%  '
$ws.Range("D30").Value = '''1.00'
$ws.Range("E30").Value = '  +0.01%  '
$ws.Range("E31").Value = '  -2.00%  '
$ws.Range("D32").Value = '''7.67'
$ws.Range("E32").Value = '  -5.15%  '
$ws.Range("D33").Value = '''22.89'
$ws.Range("E33").Value = '  -2.34%  '
$ws.Range("D34").Value = '''1.29'
$ws.Range("E34").Value = '  -3.50%  '
$ws.Range("D35").Value = '''5.37'
$ws.Range("E35").Value = '  +0.23%  '
$ws.Range("D36").Value = '''168.39'
$ws.Range("E36").Value = '  -0.98%  '
$ws.Range("E37").Value = '  -5.81%  '
$ws.Range("D38").Value = '''6.76'
$ws.Range("E38").Value = '  -4.05%  '
$ws.Range("D39").Value = '''29.20'
$ws.Range("E39").Value = '  -8.68%  '
$ws.Range("D40").Value = '3.391.76'
$ws.Range("E40").Value = '  -3.02%  '
$ws.Range("D41").Value = '''0.0752'
$ws.Range("E41").Value = '  -2.71%  '
$ws.Range("D42").Value = '''42.30'
$ws.Range("E42").Value = '  -1.24%  '
$ws.Range("E43").Value = '  -4.35%  '
$ws.Range("D44").Value = '''4.30'
$ws.Range("E44").Value = '  -1.64%  '
$ws.Range("E45").Value = '  -4.78%  '
$ws.Range("D46").Value = '''1.61'
$ws.Range("E46").Value = '  -6.37%  '
$ws.Range("D47").Value = '2.500.02'
$ws.Range("E47").Value = '  -3.19%  '
$ws.Range("D48").Value = '''22.66'
$ws.Range("E48").Value = '  +0.13%  '
$ws.Range("D49").Value = '''6.69'
$ws.Range("E49").Value = '  -3.12%  '
$ws.Range("E50").Value = '  -0.07%  '
$ws.Range("E51").Value = '  -2.56%  '
